# Econ_Data.xlsx edit: strip the manual leading-space indentation that was
# baked into column A's text and replace it with real formatting instead
# (left-aligned cell style + a wide, fixed column so the outline levels in
# the data are preserved visually through the column width rather than
# through leading spaces in the string itself).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Column A text, rows 1-50, with the leading spaces removed.
# ---------------------------------------------------------------------
$values = @(
    'Category',
    'National income',
    'Compensation of employees',
    'Wages and salaries',
    'Government',
    'Other',
    'Supplements to wages and salaries',
    'Employer contributions for employee pension and insurance funds1',
    'Employer contributions for government social insurance',
    'Proprietors'' income with IVA and CCAdj',
    'Farm',
    'Nonfarm',
    'Rental income of persons with CCAdj',
    'Corporate profits with IVA and CCAdj',
    'Taxes on corporate income',
    'Profits after tax with IVA and CCAdj',
    'Net dividends',
    'Undistributed profits with IVA and CCAdj',
    'Net interest and miscellaneous payments',
    'Taxes on production and imports',
    'Less: Subsidies2',
    'Business current transfer payments (net)',
    'To persons (net)',
    'To government (net)',
    'To the rest of the world (net)',
    'Current surplus of government enterprises2',
    'Net cash flow with IVA',
    'Undistributed profits with IVA and CCAdj',
    'Consumption of fixed capital',
    'Less: Capital transfers paid (net)',
    'Proprietors'' income with IVA and CCAdj',
    'Farm',
    'Proprietors'' income with IVA',
    'Capital consumption adjustment',
    'Nonfarm',
    'Proprietors'' income (without IVA and CCAdj)',
    'Inventory valuation adjustment',
    'Capital consumption adjustment',
    'Rental income of persons with CCAdj',
    'Rental income of persons (without CCAdj)',
    'Capital consumption adjustment',
    'Corporate profits with IVA and CCAdj',
    'Corporate profits with IVA',
    'Profits before tax (without IVA and CCAdj)',
    'Taxes on corporate income',
    'Profits after tax (without IVA and CCAdj)',
    'Net dividends',
    'Undistributed profits (without IVA and CCAdj)',
    'Inventory valuation adjustment',
    'Capital consumption adjustment'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# ---------------------------------------------------------------------
# 2. New cell style for column A: left-horizontal alignment (this is
#    what now conveys the outline depth instead of leading spaces).
# ---------------------------------------------------------------------
$ws.Range("A1:A50").HorizontalAlignment = -4131   # xlLeft

# ---------------------------------------------------------------------
# 3. Widen column A to fit the longest label, with a fixed custom width.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 72.8

# ---------------------------------------------------------------------
# 4. Restore the cursor/selection to A45, matching the saved view state.
# ---------------------------------------------------------------------
$ws.Range("A45").Select()
